# Loan RBI, Variable Instalments
# - Insert a new (blank) column before column N on the "Repayment schedule"
#   sheet, shifting the old N/O/P ("Late"/"Paid Date"/"Disbursement") columns
#   one place to the right.
# - Make "Repayment schedule" the active sheet/tab (it was "Transactions").
# - Update the selected cell on "Repayment schedule" to R7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert the new blank column at N; everything from N onward shifts right.
$ws.Columns("N").Insert()

# The new column inherits the width of the column to its left (M).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth()

# Make "Repayment schedule" the active sheet (was "Transactions").
$ws.Activate()

# Update selection to R7 on the now-active sheet.
$ws.Range("R7").Select() | Out-Null
